# Recreates the "Add files via upload" commit: a new worksheet ("Sheet1")
# is appended to the workbook containing a 26-line Splunk SPL query (one
# line per row in column A), and the view/selection state of the existing
# "set-adminauditlogconfig" sheet is updated.

$wb = $excel.ActiveWorkbook

# --- 1. Update the selection on the existing "set-adminauditlogconfig" sheet ---
$auditSheet = $wb.Worksheets.Item("set-adminauditlogconfig")
$auditSheet.Select()
$auditSheet.Range("F23").Select()

# --- 2. Add the new worksheet at the end of the workbook and rename it ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "Sheet1"

# --- 3. Populate column A with the SPL query, one line per row ---
$lines = @(
    '| makeresults',
    '| eval CreationTime="2025-06-10T16:23:42",',
    '      Id="2af83261-993a-4579-da33-08dda83b2a86",',
    '      Operation="Set-AdminAuditLogConfig",',
    '      OrganizationId="384bb877-d355-46c9-a24e-cfdad23dd148",',
    '      RecordType=1,',
    '      ResultStatus="True",',
    '      UserKey="admin@domain.onmicrosoft.com",',
    '      UserType=2,',
    '      Version=1,',
    '      Workload="Exchange",',
    '      ClientIP="127.0.0.1:15166",',
    '      ObjectId="Admin Audit Log Settings",',
    '      UserId="admin@domain.onmicrosoft.com",',
    '      AppId="",',
    '      AppPoolName="MSExchangeAdminApiNetCore",',
    '      ClientAppId="",',
    '      CorrelationID="",',
    '      ExternalAccess=false,',
    '      OrganizationName="domain.onmicrosoft.com",',
    '      OriginatingServer="BY5PR17MB3779 (15.20.8746.028)",',
    '      Parameters_Name="UnifiedAuditLogIngestionEnabled",',
    '      Parameters_Value="False",',
    '      RequestId="72f49cd1-b98e-dbe5-4481-ce8a783efd82",',
    '      SessionId="005c58b9-f35d-dd21-a262-881d4e6bcd9d"',
    '| table CreationTime Id Operation OrganizationId RecordType ResultStatus UserKey UserType Version Workload ClientIP ObjectId UserId AppId AppPoolName ClientAppId CorrelationID ExternalAccess OrganizationName OriginatingServer Parameters_Name Parameters_Value RequestId SessionId'
)

for ($i = 0; $i -lt $lines.Length; $i++) {
    $newSheet.Cells.Item($i + 1, 1).Value = $lines[$i]
}

# --- 4. Cosmetic sheet formatting matching the authored workbook ---
$newSheet.Columns.Item(1).ColumnWidth = 63.79

# --- 5. Leave the selection on the new sheet where the author left it ---
$newSheet.Range("A4").Select()
